## TC_114_Verify_Cable_capacitance.xlsx - apply "Updated test data for normal
## load, cable capacitance etc" edit to the "Devices" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")
$ws.Activate()

# --- New header cell H1 ("ColumnNumber") -----------------------------------
# Must be written AFTER the I-column strings below so the shared-string
# table ends up in the same order as the target file (Label Name,
# Exi800 - 1, Exi800 - 2, Exi800 - 3, ColumnNumber).

# --- New column I: "Label Name" header + Exi800 labels ---------------------
$ws.Range("I7").Value = "Label Name"
$ws.Range("H7").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null

$ws.Range("I8").Value = "Exi800 - 1"
$ws.Range("A8").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null

$ws.Range("I9").Value = "Exi800 - 2"
$ws.Range("A9").Copy() | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null

$ws.Range("I10").Value = "Exi800 - 3"
$ws.Range("A10").Copy() | Out-Null
$ws.Range("I10").PasteSpecial(-4122) | Out-Null

# --- New column H additions (H1 header, H2 value) ---------------------------
$ws.Range("H1").Value = "ColumnNumber"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

$ws.Range("H2").Value = 1
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row height / column widths --------------------------------------------
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Columns.Item(2).ColumnWidth = 23
$ws.Columns.Item(6).ColumnWidth = 20

# --- Selection ---------------------------------------------------------------
$ws.Range("H4").Select() | Out-Null
